$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows at the top of this subset (rows 704-705), pushing
# all existing rows (old 704-747) down to (706-749).
$ws.Range("A704:R705").EntireRow.Insert()

# --- New row 704 ---
$ws.Range("A704").Value = 8
$ws.Range("B704").Value = "Terminal La Palmera de La Serena"
$ws.Range("C704").Value = "Coquimbo"
$ws.Range("D704").Value = 44931
$ws.Range("E704").Value = 4
$ws.Range("F704").Value = 100112043
$ws.Range("G704").Value = "Pepino ensalada"
$ws.Range("H704").Value = "Sin especificar"
$ws.Range("I704").Value = "Primera"
$ws.Range("J704").Value = 700
$ws.Range("K704").Value = 12000
$ws.Range("L704").Value = 13000
$ws.Range("M704").Value = 12500
$ws.Range("N704").Value = "$/caja 60 unidades"
$ws.Range("O704").Value = "Región de Arica y Parinacota"
$ws.Range("P704").Value = 208
$ws.Range("Q704").Value = 60
$ws.Range("R704").Value = "Hortaliza"

# --- New row 705 ---
$ws.Range("A705").Value = 8
$ws.Range("B705").Value = "Terminal La Palmera de La Serena"
$ws.Range("C705").Value = "Coquimbo"
$ws.Range("D705").Value = 44931
$ws.Range("E705").Value = 4
$ws.Range("F705").Value = 100112043
$ws.Range("G705").Value = "Pepino ensalada"
$ws.Range("H705").Value = "Sin especificar"
$ws.Range("I705").Value = "Segunda"
$ws.Range("J705").Value = 400
$ws.Range("K705").Value = 10000
$ws.Range("L705").Value = 11000
$ws.Range("M705").Value = 10500
$ws.Range("N705").Value = "$/caja 80 unidades"
$ws.Range("O705").Value = "Región de Arica y Parinacota"
$ws.Range("P705").Value = 131
$ws.Range("Q705").Value = 80
$ws.Range("R705").Value = "Hortaliza"
